# Insert a new data row after the header/initial rows, at row 7.
# This pushes the existing rows 7..38 down to 8..39 (Excel keeps cell
# formatting/merges consistent automatically), then we populate the
# newly inserted row 7 with its data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(7).Insert()

$ws.Cells.Item(7, 1).Value = 7
$ws.Cells.Item(7, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(7, 3).Value = "Ñuble"
$ws.Cells.Item(7, 4).Value = 45030
$ws.Cells.Item(7, 5).Value = 16
$ws.Cells.Item(7, 6).Value = 100112043
$ws.Cells.Item(7, 7).Value = "Pepino dulce"
$ws.Cells.Item(7, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 60
$ws.Cells.Item(7, 11).Value = 15000
$ws.Cells.Item(7, 12).Value = 15000
$ws.Cells.Item(7, 13).Value = 15000
$ws.Cells.Item(7, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(7, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(7, 16).Value = 833
$ws.Cells.Item(7, 17).Value = 18
$ws.Cells.Item(7, 18).Value = "Hortaliza"
